$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@('terqui1415', '2024-06-15', 'paulo', 'r, a, f, a, e, l')
    ,@('segqua2122', '2024-06-16', 'pedro', 't, h, a, l, e, s')
    ,@('segqua2122', '2024-06-16', 'pedro', 't, h, a, l, e, s')
    ,@('terqui0910', '2024-06-18', 'julio', 'n, a, t, a, s, h, a')
    ,@('segqua0708', '2024-06-18', 'paulo', 'p, i, p, i, c, o')
    ,@('terqui1415', '2024-06-30', 'julio', 'rafael')
    ,@('terqui1415', '2024-06-30', 'julio', 'rafael')
    ,@('terqui1415', '2024-06-30', 'julio', 'rafael')
    ,@('terqui1415', '2024-06-30', 'julio', 'rafael')
    ,@('terqui1415', '2024-06-30', 'julio', 'rafael')
    ,@('terqui1415', '2024-06-26', 'julio', 'jennifer')
    ,@('terqui1314', '2024-06-26', 'julio', 'jennifer')
    ,@('terqui1415', '2024-06-24', 'pedro', 'jennifer, luiz')
    ,@('terqui1415', '2024-06-24', 'pedro', 'jennifer, luiz')
    ,@('segqua2122', '2024-06-26', 'julio', 'thales, letícia')
    ,@('segqua1112', '2024-06-14', 'julio', 'andré')
    ,@('segqua0708', '2024-06-30', 'pedro', 'manel, bernardo, pipico')
    ,@('segqua0708', '2024-06-16', 'pedro', 'Manel, Bernardo, Pipico')
    ,@('segqua2122', '2024-06-11', 'pedro', 'Thales, Amanda')
)

$startRow = 3
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    # Column B holds dates formatted as plain text (e.g. "2024-06-15"), not real
    # date serials - force text storage, then drop the residual text-format style
    # so the cell ends up with no explicit style, matching the source data.
    $ws.Cells.Item($r, 2).NumberFormat = "@"
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 2).ClearFormats()
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
